$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-08-19 Tuesday" "2025-08-20 Wednesday"

Replace-Text "478÷3=159, 1" "480÷4=120, 0"
Replace-Text "720÷8=90, 0" "655÷4=163, 3"
Replace-Text "844÷3=281, 1" "123÷6=20, 3"
Replace-Text "928÷9=103, 1" "496÷4=124, 0"
Replace-Text "396÷3=132, 0" "348÷9=38, 6"

Replace-Text "735÷9=81, 6" "850÷5=170, 0"
Replace-Text "599÷5=119, 4" "867÷5=173, 2"
Replace-Text "825÷7=117, 6" "722÷2=361, 0"
Replace-Text "761÷5=152, 1" "448÷6=74, 4"
Replace-Text "799÷8=99, 7" "963÷3=321, 0"

Replace-Text "316÷9=35, 1" "958÷7=136, 6"
Replace-Text "863÷9=95, 8" "330÷4=82, 2"
Replace-Text "234÷4=58, 2" "344÷6=57, 2"
Replace-Text "627÷9=69, 6" "776÷4=194, 0"
Replace-Text "979÷7=139, 6" "630÷9=70, 0"

Replace-Text "761÷4=190, 1" "911÷4=227, 3"
Replace-Text "419÷2=209, 1" "264÷9=29, 3"
Replace-Text "514÷5=102, 4" "832÷6=138, 4"
Replace-Text "593÷9=65, 8" "464÷9=51, 5"
Replace-Text "676÷2=338, 0" "723÷3=241, 0"

Replace-Text "395÷7=56, 3" "243÷4=60, 3"
Replace-Text "480÷5=96, 0" "630÷4=157, 2"
Replace-Text "802÷5=160, 2" "138÷2=69, 0"
Replace-Text "241÷6=40, 1" "670÷4=167, 2"
Replace-Text "615÷6=102, 3" "232÷7=33, 1"
